# Apply the "Anonimyzed fedcore" update:
#  - rename "fedcore" header labels to "approach"
#  - add a thin top+bottom border (C1/F1) and a thin top+bottom+right border
#    (D1/G1) to the merged-header spacer cells on both sheets
#  - drop the stray empty inline-string cell at G5 on the computational sheet

$wb = $excel.ActiveWorkbook

$quality = $wb.Worksheets.Item("quality_comparison")
$comp    = $wb.Worksheets.Item("computational_comparison")

function Set-TopBottomBorder($range) {
    $range.Style = "Normal"
    $range.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $range.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
}

function Set-TopBottomRightBorder($range) {
    $range.Style = "Normal"
    $range.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $range.Borders.Item(10).LineStyle = 1  # xlEdgeRight
    $range.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
}

# --- quality_comparison sheet ---
Set-TopBottomBorder      $quality.Range("C1")
Set-TopBottomRightBorder $quality.Range("D1")
$quality.Range("C2").Value = "approach"

# --- computational_comparison sheet ---
Set-TopBottomBorder      $comp.Range("C1")
Set-TopBottomRightBorder $comp.Range("D1")
Set-TopBottomBorder      $comp.Range("F1")
Set-TopBottomRightBorder $comp.Range("G1")
$comp.Range("C2").Value = "approach"
$comp.Range("F2").Value = "approach"

# remove the stray empty inline-string cell
$comp.Range("G5").ClearContents()
